# AddToCart test and page with generic xpaths
# Splits the "Product Details" sheet's rows 3 & 4 (Computers/Desktops/HP
# TOUCHSMART and Apparels/PURE COTTON T-SHIRT) out into a brand-new
# "Sheet1" worksheet (placed between "Product Details" and
# "SignIn_Data_Negative"), shifted two columns to the right (C:H) and
# sitting at rows 5-6. Also renames "Sports" -> "Sport" and
# "Apparels" -> "Apparel", and re-types the product-number cell for the
# first product as a quote-prefixed text value "133".

$wb = $excel.ActiveWorkbook
$wsProduct = $wb.Worksheets.Item("Product Details")
$wsNegative = $wb.Worksheets.Item("SignIn_Data_Negative")

# --- 1. Rename "Sports" -> "Sport" first, so the shared-string table ends
#        up with the same append order the source workbook has (Sport,
#        then Apparel, then "133") once the now-unused old strings are
#        compacted away on save. ------------------------------------------
$wsProduct.Range("A2").Value = "Sport"

# --- 2. Insert the new worksheet between "Product Details" and
#        "SignIn_Data_Negative" ------------------------------------------
$wsNew = $wb.Worksheets.Add($null, $wsProduct)
$wsNew.Name = "Sheet1"

# --- 3. Populate the new sheet with the data that used to live in rows
#        3-4 of "Product Details" (columns A:F), now living at rows 5-6,
#        columns C:H -------------------------------------------------------
$wsNew.Range("C5").Value = "Computers"
$wsNew.Range("D5").Value = "Desktops"
$wsNew.Range("E5").Value = "HP TOUCHSMART"
$wsNew.Range("F5").Value = 31
$wsNew.Range("G5").Value = 799.99
$wsNew.Range("H5").Value = 1

$wsNew.Range("C6").Value = "Apparel"
$wsNew.Range("E6").Value = "PURE COTTON T-SHIRT"
$wsNew.Range("F6").Value = 47
$wsNew.Range("G6").Value = 13.99
$wsNew.Range("H6").Value = 2

# Match the source formatting (numFmtId 49 "@" / style index 2 on the
# original sheet) on every cell in the block, including the blank D6.
$wsNew.Range("C5:H6").NumberFormat = "@"

$wsNew.Range("C5:H6").Select()

# --- 4. Remove the now-duplicated rows 3 & 4 from "Product Details" -------
$wsProduct.Rows.Item(4).Delete()
$wsProduct.Rows.Item(3).Delete()

# --- 5. Re-type the product number for row 2 as quote-prefixed text ------
$wsProduct.Range("D2").NumberFormat = "@"
$wsProduct.Range("D2").Value = "'133"

$wsProduct.Range("D3").Select()

Write-Output "done"
